# collect from PONDR checked
# Update the "human_order" (column G) and "human_disorder" (column H) counts
# on Sheet1 with refreshed values collected from PONDR.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2"  = 479203
    "G3"  = 175643
    "G4"  = 330080
    "H4"  = 207956
    "G5"  = 399940
    "H5"  = 406608
    "G6"  = 324013
    "G7"  = 431350
    "G8"  = 191475
    "G9"  = 378838
    "H9"  = 113800
    "G10" = 360043
    "H10" = 290938
    "G11" = 771505
    "H11" = 359909
    "G12" = 159982
    "G13" = 274701
    "H13" = 132976
    "G14" = 320476
    "G15" = 295531
    "H15" = 246048
    "G16" = 355342
    "G17" = 477217
    "H17" = 468701
    "G18" = 370569
    "H18" = 237125
    "G19" = 478562
    "H19" = 198763
    "G20" = 113612
    "G21" = 238289
    "H21" = 64256
    "G22" = 4429825
    "H22" = 6926371
    "G23" = 11356196
    "H23" = 11356196
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
